$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the data refresh: updated prices / 1h-volume percentages, and the
# two rank swaps (Chainlink <-> WrappedliquidstakedEther2.0 at rows 13/14,
# RenderToken <-> RocketPoolETH at rows 48/49).
$ws.Range("D2").Value = '36.809.19'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '2.116.80'
$ws.Range("E3").Value = '  +10.10%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'256.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.68%  '
$ws.Range("D6").Value = "'0.669"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.67%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'46.54"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +6.08%  '
$ws.Range("D9").Value = "'62.25"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.68%  '
$ws.Range("D10").Value = "'0.373"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = "'0.0744"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.424.62'
$ws.Range("E13").Value = '  +10.17%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'14.59"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = "'0.853"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +6.81%  '
$ws.Range("D16").Value = '2.114.71'
$ws.Range("E16").Value = '  +9.91%  '
$ws.Range("D17").Value = "'5.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").Value = '36.786.83'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").Value = "'74.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = "'13.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.76%  '
$ws.Range("D22").Value = "'242.22"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.46%  '
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").Value = '  -8.07%  '
$ws.Range("D26").Value = "'173.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.99%  '
$ws.Range("D27").Value = "'21.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +13.43%  '
$ws.Range("D28").Value = "'9.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.86%  '
$ws.Range("D29").Value = "'2.05"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -9.39%  '
$ws.Range("E30").Value = '  -3.98%  '
$ws.Range("D31").Value = "'22.64"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +50.36%  '
$ws.Range("D32").Value = "'4.58"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("D33").Value = "'0.0961"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +13.84%  '
$ws.Range("D34").Value = "'0.0606"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").Value = "'2.42"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +19.68%  '
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = "'4.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("D39").Value = "'0.921"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.60%  '
$ws.Range("E40").Value = '  -8.27%  '
$ws.Range("D41").Value = "'1.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +8.01%  '
$ws.Range("D42").Value = "'0.0224"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("D43").Value = "'99.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.35%  '
$ws.Range("D44").Value = "'2.83"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +17.15%  '
$ws.Range("D45").Value = "'16.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.25%  '
$ws.Range("D46").Value = '1.365.05'
$ws.Range("E46").Value = '  +1.52%  '
$ws.Range("D47").Value = "'0.0839"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.63%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.316.32'
$ws.Range("E48").Value = '  +10.28%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'2.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.17%  '
$ws.Range("E50").Value = '  +6.96%  '
$ws.Range("D51").Value = "'2.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.67%  '
